$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / header rich-text updates ---
# A8: "Volume 32   Number  34" -> "...35" (change run text "34"->"35")
$ws.Range("A8").Characters(21,2).Text = "35"
# C9: date range update within rich text runs
$ws.Range("C9").Characters(27,9).Text = "8/25/2025"
$ws.Range("C9").Characters(47,9).Text = "8/31/2025"

# --- Crime data table updates (rows 15-28) ---
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = 0
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 1
$ws.Range("G15").NumberFormat = '#,##0'
$ws.Range("H15").Value = 200
$ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I15").Value = 23
$ws.Range("J15").Value = 17
$ws.Range("K15").Value = 35.294117647058
$ws.Range("L15").Value = 53.333333333333
$ws.Range("M15").Value = 64.285714285714
$ws.Range("N15").Value = 27.777777777777
$ws.Range("C16").Value = 8
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -22.222222222222
$ws.Range("I16").Value = 96
$ws.Range("J16").Value = 141
$ws.Range("K16").Value = -31.914893617021
$ws.Range("L16").Value = -34.246575342465
$ws.Range("M16").Value = -46.961325966850
$ws.Range("N16").Value = -84.565916398713
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 16
$ws.Range("H17").Value = -27.272727272727
$ws.Range("I17").Value = 208
$ws.Range("J17").Value = 202
$ws.Range("K17").Value = 2.970297029702
$ws.Range("L17").Value = 13.661202185792
$ws.Range("M17").Value = 35.947712418300
$ws.Range("N17").Value = 0.970873786407
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 14.285714285714
$ws.Range("I18").Value = 110
$ws.Range("J18").Value = 151
$ws.Range("K18").Value = -27.152317880794
$ws.Range("L18").Value = -14.728682170542
$ws.Range("M18").Value = -63.934426229508
$ws.Range("N18").Value = -92.052023121387
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 45.454545454545
$ws.Range("F19").Value = 55
$ws.Range("H19").Value = 25
$ws.Range("I19").Value = 436
$ws.Range("J19").Value = 449
$ws.Range("K19").Value = -2.895322939866
$ws.Range("L19").Value = -5.627705627705
$ws.Range("M19").Value = 59.124087591240
$ws.Range("N19").Value = -0.909090909090
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 15
$ws.Range("E20").Value = -40
$ws.Range("F20").Value = 41
$ws.Range("H20").Value = -2.380952380952
$ws.Range("I20").Value = 243
$ws.Range("J20").Value = 288
$ws.Range("K20").Value = -15.625
$ws.Range("L20").Value = -0.409836065573
$ws.Range("M20").Value = -6.177606177606
$ws.Range("N20").Value = -89.958677685950
$ws.Range("C21").Value = 42
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = 16.666666666666
$ws.Range("F21").Value = 145
$ws.Range("G21").Value = 141
$ws.Range("H21").Value = 2.836879432624
$ws.Range("I21").Value = 1120
$ws.Range("J21").Value = 1249
$ws.Range("K21").Value = -10.328262610088
$ws.Range("L21").Value = -5.245346869712
$ws.Range("M21").Value = -5.644481887110
$ws.Range("N21").Value = -78.065021543282
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("F22").Value = 2
$ws.Range("I22").Value = 12
$ws.Range("K22").Value = -7.692307692307
$ws.Range("L22").Value = -29.411764705882
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 38
$ws.Range("E24").Value = -18.421052631578
$ws.Range("F24").Value = 114
$ws.Range("G24").Value = 115
$ws.Range("H24").Value = -0.869565217391
$ws.Range("I24").Value = 832
$ws.Range("J24").Value = 977
$ws.Range("K24").Value = -14.841351074718
$ws.Range("L24").Value = -6.830907054871
$ws.Range("M24").Value = 15.395284327323
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 22
$ws.Range("E25").Value = -59.090909090909
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 61
$ws.Range("H25").Value = -32.786885245901
$ws.Range("I25").Value = 318
$ws.Range("J25").Value = 384
$ws.Range("K25").Value = -17.1875
$ws.Range("L25").Value = 1.923076923076
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 12
$ws.Range("E26").Value = 16.666666666666
$ws.Range("F26").Value = 54
$ws.Range("G26").Value = 41
$ws.Range("H26").Value = 31.707317073170
$ws.Range("I26").Value = 373
$ws.Range("J26").Value = 400
$ws.Range("K26").Value = -6.75
$ws.Range("L26").Value = 16.199376947040
$ws.Range("M26").Value = -24.949698189134
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = 0
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 1
$ws.Range("G27").NumberFormat = '#,##0'
$ws.Range("H27").Value = 300
$ws.Range("H27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I27").Value = 25
$ws.Range("J27").Value = 23
$ws.Range("K27").Value = 8.695652173913
$ws.Range("L27").Value = 4.166666666666
$ws.Range("C23").Copy($ws.Range("C28"))
$ws.Range("C23").Copy($ws.Range("D28"))
$ws.Range("E23").Copy($ws.Range("E28"))
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = -50

